# Re-randomized trial data for subject 11 / living_rooms categorization block 2.
# Each row's trial_total (F), stimulus (L), the conceptual/perceptual/typicality
# scores (M:O), sample size n (P) and the three proportion columns (Q:S) are
# refreshed in place; category/cond_cat/correct_answer (H/I/K) only change for
# rows whose stimulus switched between target and distractor categories.
# This mirrors 'elaborated sanity checks. increased the proportion of new
# images in memory task.'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 135
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_wgkqa.png'
$arr[0,7] = 87.25581395348837
$arr[0,8] = 71.13953488372093
$arr[0,9] = 79.19767441860465
$arr[0,10] = 43
$arr[0,11] = 10
$arr[0,12] = 10
$arr[0,13] = 10
$ws.Range("F2:S2").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 136
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_eh0no.png'
$arr[0,7] = 53.66666666666666
$arr[0,8] = 36.02564102564103
$arr[0,9] = 44.84615384615385
$arr[0,10] = 39
$arr[0,11] = 3
$arr[0,12] = 3
$arr[0,13] = 3
$ws.Range("F3:S3").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 137
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'kitchens'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_pt3d7.png'
$arr[0,7] = 65.08571428571429
$arr[0,8] = 44.65714285714286
$arr[0,9] = 54.87142857142857
$arr[0,10] = 35
$arr[0,11] = 4
$arr[0,12] = 4
$arr[0,13] = 4
$ws.Range("F4:S4").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 138
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_di6f0.png'
$arr[0,7] = 94.04347826086956
$arr[0,8] = 83.34782608695652
$arr[0,9] = 88.69565217391303
$arr[0,10] = 46
$arr[0,11] = 10
$arr[0,12] = 10
$arr[0,13] = 10
$ws.Range("F5:S5").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 139
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_16kib.png'
$arr[0,7] = 80.97727272727273
$arr[0,8] = 61.11363636363637
$arr[0,9] = 71.04545454545455
$arr[0,10] = 44
$arr[0,11] = 8
$arr[0,12] = 8
$arr[0,13] = 8
$ws.Range("F6:S6").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 140
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'kitchens'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_kugyw.png'
$arr[0,7] = 74.25
$arr[0,8] = 54.10714285714285
$arr[0,9] = 64.17857142857143
$arr[0,10] = 28
$arr[0,11] = 6
$arr[0,12] = 6
$arr[0,13] = 6
$ws.Range("F7:S7").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 141
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_cehin.png'
$arr[0,7] = 78.86363636363636
$arr[0,8] = 60.02272727272727
$arr[0,9] = 69.44318181818181
$arr[0,10] = 44
$arr[0,11] = 7
$arr[0,12] = 7
$arr[0,13] = 7
$ws.Range("F8:S8").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 142
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_4o8l0.png'
$arr[0,7] = 46.02173913043478
$arr[0,8] = 31.45652173913043
$arr[0,9] = 38.73913043478261
$arr[0,10] = 46
$arr[0,11] = 3
$arr[0,12] = 3
$arr[0,13] = 3
$ws.Range("F9:S9").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 143
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_xu1p3.png'
$arr[0,7] = 75.27659574468085
$arr[0,8] = 56.68085106382978
$arr[0,9] = 65.97872340425532
$arr[0,10] = 47
$arr[0,11] = 7
$arr[0,12] = 7
$arr[0,13] = 7
$ws.Range("F10:S10").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 144
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_6a0hu.png'
$arr[0,7] = 61.275
$arr[0,8] = 42.025
$arr[0,9] = 51.65
$arr[0,10] = 40
$arr[0,11] = 4
$arr[0,12] = 4
$arr[0,13] = 4
$ws.Range("F11:S11").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 145
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_pey7u.png'
$arr[0,7] = 30.34883720930232
$arr[0,8] = 20.34883720930232
$arr[0,9] = 25.34883720930232
$arr[0,10] = 43
$arr[0,11] = 1
$arr[0,12] = 2
$arr[0,13] = 2
$ws.Range("F12:S12").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 146
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_xy930.png'
$arr[0,7] = 70.5952380952381
$arr[0,8] = 49.47619047619047
$arr[0,9] = 60.03571428571429
$arr[0,10] = 42
$arr[0,11] = 6
$arr[0,12] = 6
$arr[0,13] = 6
$ws.Range("F13:S13").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 147
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_abobq.png'
$arr[0,7] = 75.1842105263158
$arr[0,8] = 54.13157894736842
$arr[0,9] = 64.65789473684211
$arr[0,10] = 38
$arr[0,11] = 6
$arr[0,12] = 6
$arr[0,13] = 6
$ws.Range("F14:S14").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 148
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'bedrooms'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_ys3qz.png'
$arr[0,7] = 46.79545454545455
$arr[0,8] = 31.20454545454545
$arr[0,9] = 39
$arr[0,10] = 44
$arr[0,11] = 2
$arr[0,12] = 2
$arr[0,13] = 2
$ws.Range("F15:S15").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 149
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'bedrooms'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_fea1z.png'
$arr[0,7] = 79.45945945945945
$arr[0,8] = 56.24324324324324
$arr[0,9] = 67.85135135135135
$arr[0,10] = 37
$arr[0,11] = 7
$arr[0,12] = 7
$arr[0,13] = 7
$ws.Range("F16:S16").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 150
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'kitchens'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_cxpff.png'
$arr[0,7] = 74.92307692307692
$arr[0,8] = 53.28205128205128
$arr[0,9] = 64.1025641025641
$arr[0,10] = 39
$arr[0,11] = 6
$arr[0,12] = 6
$arr[0,13] = 6
$ws.Range("F17:S17").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 151
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_w8yhd.png'
$arr[0,7] = 55.74418604651163
$arr[0,8] = 38.90697674418605
$arr[0,9] = 47.32558139534883
$arr[0,10] = 43
$arr[0,11] = 4
$arr[0,12] = 4
$arr[0,13] = 4
$ws.Range("F18:S18").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 152
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_xbtev.png'
$arr[0,7] = 13.68181818181818
$arr[0,8] = 8.568181818181818
$arr[0,9] = 11.125
$arr[0,10] = 44
$arr[0,11] = 1
$arr[0,12] = 1
$arr[0,13] = 1
$ws.Range("F19:S19").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 153
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_6zz63.png'
$arr[0,7] = 87.66666666666667
$arr[0,8] = 70.6
$arr[0,9] = 79.13333333333333
$arr[0,10] = 45
$arr[0,11] = 9
$arr[0,12] = 10
$arr[0,13] = 10
$ws.Range("F20:S20").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 154
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_amsgw.png'
$arr[0,7] = 86.08510638297872
$arr[0,8] = 65.95744680851064
$arr[0,9] = 76.02127659574468
$arr[0,10] = 47
$arr[0,11] = 9
$arr[0,12] = 9
$arr[0,13] = 9
$ws.Range("F21:S21").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 155
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_0kqc0.png'
$arr[0,7] = 43.74468085106383
$arr[0,8] = 27.14893617021277
$arr[0,9] = 35.4468085106383
$arr[0,10] = 47
$arr[0,11] = 2
$arr[0,12] = 2
$arr[0,13] = 2
$ws.Range("F22:S22").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 156
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_bj99b.png'
$arr[0,7] = 82.79069767441861
$arr[0,8] = 65.46511627906976
$arr[0,9] = 74.12790697674419
$arr[0,10] = 43
$arr[0,11] = 8
$arr[0,12] = 8
$arr[0,13] = 8
$ws.Range("F23:S23").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 157
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_kost0.png'
$arr[0,7] = 63.09090909090909
$arr[0,8] = 42.77272727272727
$arr[0,9] = 52.93181818181819
$arr[0,10] = 44
$arr[0,11] = 5
$arr[0,12] = 5
$arr[0,13] = 5
$ws.Range("F24:S24").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 158
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_wz6x5.png'
$arr[0,7] = 68.3695652173913
$arr[0,8] = 48.47826086956522
$arr[0,9] = 58.42391304347826
$arr[0,10] = 46
$arr[0,11] = 5
$arr[0,12] = 5
$arr[0,13] = 5
$ws.Range("F25:S25").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 159
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'bedrooms'
$arr[0,3] = 'distractor'
$arr[0,4] = $null
$arr[0,5] = 'f'
$arr[0,6] = 'stimuli/img_twj5p.png'
$arr[0,7] = 67.71739130434783
$arr[0,8] = 42.08695652173913
$arr[0,9] = 54.90217391304348
$arr[0,10] = 46
$arr[0,11] = 4
$arr[0,12] = 4
$arr[0,13] = 4
$ws.Range("F26:S26").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 160
$arr[0,1] = 'living_rooms'
$arr[0,2] = 'living_rooms'
$arr[0,3] = 'target'
$arr[0,4] = $null
$arr[0,5] = 'j'
$arr[0,6] = 'stimuli/img_bbs77.png'
$arr[0,7] = 31.64444444444445
$arr[0,8] = 21.26666666666667
$arr[0,9] = 26.45555555555556
$arr[0,10] = 45
$arr[0,11] = 2
$arr[0,12] = 2
$arr[0,13] = 2
$ws.Range("F27:S27").Value = $arr
